$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '72.654.88'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.643.15'
$ws.Range('E3').Value = '  -1.49%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.77'
$ws.Range('E5').Value = '  -2.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.15'
$ws.Range('E6').Value = '  -0.67%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.521'
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.172'
$ws.Range('E9').Value = '  +1.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '2.642.41'
$ws.Range('E10').Value = '  -1.48%  '
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.360'
$ws.Range('E12').Value = '  +1.43%  '
$ws.Range('E13').Value = '  -1.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.124.05'
$ws.Range('E14').Value = '  -1.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000187'
$ws.Range('E15').Value = '  +0.64%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '72.500.11'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.81'
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.639.60'
$ws.Range('E18').Value = '  -1.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.09'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '376.91'
$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.88'
$ws.Range('E21').Value = '  -0.29%  '
$ws.Range('E22').Value = '  -1.47%  '
$ws.Range('E23').Value = '  -0.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.52'
$ws.Range('E24').Value = '  -0.95%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.25'
$ws.Range('E26').Value = '  -2.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.53'
$ws.Range('E27').Value = '  -2.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.777.90'
$ws.Range('E28').Value = '  -1.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0953'
$ws.Range('E30').Value = '  +1.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.99'
$ws.Range('E31').Value = '  -1.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '493.30'
$ws.Range('E32').Value = '  -3.63%  '
$ws.Range('E33').Value = '  +0.78%  '
$ws.Range('E34').Value = '  -1.21%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '161.63'
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('E37').Value = '  +6.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.21'
$ws.Range('E38').Value = '  -1.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.88'
$ws.Range('E39').Value = '  -1.22%  '
$ws.Range('E40').Value = '  -1.22%  '
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('E42').Value = '  -4.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.58'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('E44').Value = '  -2.28%  '
$ws.Range('E45').Value = '  -1.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.07'
$ws.Range('E46').Value = '  -0.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '150.90'
$ws.Range('E47').Value = '  -1.82%  '
$ws.Range('E48').Value = '  -1.53%  '
$ws.Range('E49').Value = '  -2.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.68'
$ws.Range('E50').Value = '  -3.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.610'
$ws.Range('E51').Value = '  +0.77%  '
